# Investments.xlsx - Table-5.1: fill in Sl.no 8 and 9 answers (rows 13 & 14)
# and update the active selection, per the commit:
# "Solution and results for sl.no 8 and 9 in table 5.1 and updated table 5.1 resutls"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table-5.1")
$ws.Activate()

# Row 13 (Sl.no 9): top sector count-wise -> company with highest investment
$ws.Range("C13").Value = "/organization/social-finance"
$ws.Range("D13").Value = "/organization/oneweb"
$ws.Range("E13").Value = "/organization/flipkart"

# Row 14 (Sl.no 10): second best sector count-wise -> company with highest investment
$ws.Range("C14").Value = "/organization/freescale"
$ws.Range("D14").Value = "/organization/powa-technologies"
$ws.Range("E14").Value = "/organization/shopclues-com"

# Move the live selection to reflect where the author ended up after entering
# the results (matches the saved <selection activeCell="E20" sqref="E20"/>)
[void]$ws.Range("E20").Select()
